$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Long, verbatim text blocks. Single-quoted here-strings => no $var expansion.
# ---------------------------------------------------------------------------
$bigPrayer = @'
*1. Покаяние*^
^
Отче наш, Отец Небесный,^
прошу Тебя простить^
все мои прегрешения,^
вольные или невольные.^
^
Аминь!^
^
*2. Новая Молитва*^
^
ОТЧЕ НАШ, ОТЕЦ НЕБЕСНЫЙ,^
Я ПРИНИМАЮ ВЕРУ ТВОЮ,^
ОНА ЕСТЬ МОЙ ПУТЬ.^
^
Я ПРИНИМАЮ КАНОНЫ ТВОИ ВЕЧНЫЕ,^
С ЛЮБОВЬЮ К ТЕБЕ И ДЕЛАМ ТВОИМ,^
ПОДТВЕРЖДАЯ СВОЕЙ ЖИЗНЬЮ^
ВЕРНОСТЬ ТЕБЕ.^
^
ГОСПОДИ, ПРОШУ ДАТЬ МНЕ НАДЕЖДУ^
НА СПАСЕНИЕ ДУШИ МОЕЙ,^
И ДАРОВАТЬ МУДРОСТЬ ТВОЮ^
ДЛЯ ЖИЗНИ МОЕЙ ЗДЕСЬ,^
НА ПЛАНЕТЕ СВЯТАЯ РУСЬ И В ВЕЧНОСТИ.^
^
ПУСТЬ СВЯТА БУДЕТ УВЕРЕННОСТЬ МОЯ,^
ЧТО ТЫ ЕСМЬ!^
^
Господи, я Люблю Тебя, Благодарю Тебя и Уповаю на Милость Твою! Аминь!
'@

$longMsg = @'
*Мы, Единый Народ России, проявляя Право Свободной Воли, в соответствии с Договором между
Создателем и Россией, запускаем Импульс Энергии Духовной волны, наполненной Равенством и
Любовью, на разрушение намерения Мировой Тьмы выстроить Мировую медицину в качестве принудительной меры по управлению здоровьем человека во вред самому человеку, и устранить участие Российской медицины во всех международных медицинских преступных организациях!*
'@

$singleExc = @'
*Отче наш, Отец Небесный! Волею Создателя, Пророка и Народа Пространство Святая Русь ЕСМЬ Равенство и Любовь Навечно! Да будет Свет Истины!*
'@

# ---------------------------------------------------------------------------
# Preserve the two special cell formats that already exist in the sheet
#   - C2  -> date-number-format style (used for the empty "date" cells)
#   - A4  -> wrap-text style (used for the long message / trigger cells)
# by stashing copies of them far outside the block we are about to rebuild,
# so Excel keeps reusing the SAME style index instead of registering new,
# duplicate ones when we paste the formats back in.
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("Z101").PasteSpecial(-4122)

# Wipe out the whole existing schedule (values + styles + row heights).
$ws.Rows("1:25").Delete()

# The donor cells shifted up along with everything below row 25 (by 25 rows).
$dateDonor = $ws.Range("Z75")
$wrapDonor = $ws.Range("Z76")

function Set-DateStyle($addr) {
    $dateDonor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
function Set-WrapStyle($addr) {
    $wrapDonor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
# ---------------------------------------------------------------------------
# Row 1 - header
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Текст:"
$ws.Range("B1").Value = "Время [мск]:"
$ws.Range("C1").Value = "Дата [мск]:"
$ws.Range("D1").Value = "Тип:"
$ws.Range("E1").Value = "Триггеры:"

# ---------------------------------------------------------------------------
# Rows 2-7: three (daily-prayer, time-slot) pairs, repeated per run.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = $bigPrayer
$ws.Range("B2").Value = "02:55 - 02:59"
$ws.Range("D2").Value = "ежедневный посыл"
Set-WrapStyle("A2")
Set-DateStyle("C2")
$ws.Rows("2").RowHeight = 409.5

$ws.Range("A3").Value = $singleExc
$ws.Range("B3").Value = "03:00 - 03:04"
$ws.Range("D3").Value = "ежедневный посыл"

$ws.Range("A4").Value = $bigPrayer
$ws.Range("B4").Value = "10:55 - 10:59"
$ws.Range("D4").Value = "ежедневный посыл"
Set-WrapStyle("A4")
Set-DateStyle("C4")
$ws.Rows("4").RowHeight = 79.5

$ws.Range("A5").Value = $singleExc
$ws.Range("B5").Value = "11:00 - 11:04"
$ws.Range("D5").Value = "ежедневный посыл"

$ws.Range("A6").Value = $bigPrayer
$ws.Range("B6").Value = "18:55 - 18:59"
$ws.Range("D6").Value = "ежедневный посыл"
Set-WrapStyle("A6")
Set-DateStyle("C6")
$ws.Rows("6").RowHeight = 409.5

$ws.Range("A7").Value = $singleExc
$ws.Range("B7").Value = "19:00 - 19:04"
$ws.Range("D7").Value = "ежедневный посыл"

# ---------------------------------------------------------------------------
# Row 25: the hourly ("часовой посыл") message, keeps its own ["08","13","26"]
# trigger-hours list in column C.
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = $longMsg
$ws.Range("B25").Value = "11:55 - 11:59"
$ws.Range("C25").Value = '["08", "13", "26"]'
$ws.Range("D25").Value = "часовой посыл"
Set-WrapStyle("A25")
Set-DateStyle("C25")
$ws.Rows("25").RowHeight = 75

# ---------------------------------------------------------------------------
# Rows 36-38: trailing, mostly-empty styled rows left over from editing.
# ---------------------------------------------------------------------------
Set-DateStyle("C36")
Set-DateStyle("C37")
Set-WrapStyle("A38")
Set-DateStyle("C38")

# Clean up the donor cells used for style-copying.
$ws.Range("Z75:Z76").Delete()
# ---------------------------------------------------------------------------
# Restore the sheet view: scrolled so row 13 is at the top, A23 selected.
# (ScrollRow/ScrollColumn are best-effort here; Select() drives the saved
# <selection> which is the part that is reliably round-tripped.)
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A23").Select()
